# Updated final client report generator with new expression overlays
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New overlay descriptions, keyed by Expression Overlay (column E) name
$descriptions = @{
    'Grounded Executor' = 'You live out your calling as a Grounded Executor — someone who thrives on structure, discipline, and dependable contribution. This means you turn ideas into reality with consistency and excellence. Others experience you as reliable, practical, and trustworthy. When aligned, your presence builds confidence and traction; when misaligned, it may become rigid, perfectionistic, or resistant to change.'
    'Relational Nurturer' = 'You live out your calling as a Relational Nurturer — someone who creates safety, trust, and care in your relationships. This means you tend to serve through empathy, attentiveness, and emotional support. Others experience you as a grounded and faithful presence. When aligned, your style fosters healing and belonging; when misaligned, it may show up as over-identification with others’ emotions or exhaustion.'
    'Reflective Thinker' = 'You live out your calling as a Reflective Thinker — someone who processes deeply and distills clarity from complexity. This means you need time, space, and solitude to develop meaningful insight. Others experience you as thoughtful, wise, and discerning. When aligned, your style reveals truth and direction; when misaligned, it may appear withdrawn, overly analytical, or disconnected.'
    'Joyful Connector' = 'You live out your calling as a Joyful Connector — someone who energizes others and infuses environments with warmth. This means you share your purpose through presence, affirmation, and community. Others experience you as uplifting, lively, and contagious in your enthusiasm. When aligned, your style brings life and encouragement; when misaligned, it may show up as scattered focus or avoidance of depth.'
    'Visionary Strategist' = 'You live out your calling as a Visionary Strategist — someone who casts bold vision and moves others toward meaningful outcomes. This means you tend to lead with clarity, pursue long-term goals, and inspire aligned action. Others experience you as forward-focused and mission-driven. When aligned, your style brings direction and momentum; when misaligned, it may manifest as over-control or burnout.'
    'Creative Builder' = 'You live out your calling as a Creative Builder — someone who brings truth and beauty to life through tangible expression. This means you are often creating, shaping, or designing things that reflect deeper values. Others experience you as inventive and inspirational. When aligned, your work invites others into wonder and insight; when misaligned, you may feel misunderstood, unproductive, or overly idealistic.'
}

# Column F (Overlay Description) is a pure function of column E (Expression Overlay):
# replace every Overlay Description cell with the new text for its overlay.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $overlay = $ws.Cells.Item($r, 5).Value()
    $ws.Cells.Item($r, 6).Value = $descriptions[$overlay]
}
